# Apply updated crypto price/volume figures (Sun Oct  1 05:53:33 UTC 2023 data refresh).
# Column D ("Price") and column E ("Volume(1h)") values are stored as plain text in the
# sheet, so numeric-looking prices are written with a leading apostrophe to force Excel to
# keep them as text (matching the original inline-string cells) instead of silently
# converting them to numbers (which would, e.g., turn "9.40" into 9.4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.103.25'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '1.681.87'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''215.17'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  +2.04%  '
$ws.Range("D9").Value = '''21.33'
$ws.Range("E9").Value = '  +5.71%  '
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").Value = '''0.0885'
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").Value = '1.917.78'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '1.679.39'
$ws.Range("E13").Value = '  -1.56%  '
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").Value = '''66.21'
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").Value = '27.096.52'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '''238.27'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").Value = '''8.13'
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("D20").Value = '0.0₃0749'
$ws.Range("E20").Value = '  +2.33%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  +1.49%  '
$ws.Range("D23").Value = '''9.40'
$ws.Range("E23").Value = '  +2.61%  '
$ws.Range("E24").Value = '  -2.44%  '
$ws.Range("D25").Value = '''146.87'
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("E26").Value = '  +0.81%  '
$ws.Range("E27").Value = '  +2.21%  '
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").Value = '1.557.48'
$ws.Range("E32").Value = '  +5.58%  '
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("E34").Value = '  +1.92%  '
$ws.Range("E35").Value = '  +2.47%  '
$ws.Range("D36").Value = '''0.604'
$ws.Range("E36").Value = '  +4.66%  '
$ws.Range("D37").Value = '''0.937'
$ws.Range("E37").Value = '  +4.64%  '
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("E39").Value = '  +2.37%  '
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").Value = '''68.78'
$ws.Range("E41").Value = '  +3.07%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").Value = '''5.65'
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("E44").Value = '  -1.76%  '
$ws.Range("D45").Value = '1.826.57'
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").Value = '''0.782'
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").Value = '''90.86'
$ws.Range("E47").Value = '  +0.51%  '
$ws.Range("E48").Value = '  +2.96%  '
$ws.Range("E49").Value = '  +0.47%  '
$ws.Range("E50").Value = '  +3.61%  '
$ws.Range("D51").Value = '''8.04'
$ws.Range("E51").Value = '  +4.64%  '
